$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices and 1h volume changes).
# Cells that hold purely numeric-looking price strings are forced to
# Text format first so Excel keeps them as literal strings (matching
# the source data which stores them as text, e.g. "8.00", "1.00").

$ws.Range("D2").Value = "63.238.91"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "2.567.81"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.29"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.14"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.77"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.45"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").Value = "3.028.56"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "63.103.98"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000154"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "2.605.50"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.15"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.67"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.39"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.86"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -4.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.09"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "556.51"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.00"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.162"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "0.0₃0856"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "165.29"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.413"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.38"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "165.33"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.63"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.96"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0591"
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.61"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.05"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.624"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0248"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0960"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.92"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "0.0₆0229"
$ws.Range("E51").Value = "  +15.22%  "
